$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ManchesterTest")

$ws.Range("A7").Value = "ex:TEST_DECIMAL"
$ws.Range("B7").Value = "ex:literal_qtt_mg some xsd:decimal[>= 2.5 , <= 3]"
